# Update "想去人数" (F column) figures on the 展览 (Exhibitions), 演出 (Shows)
# and 全部类型 (All types) sheets to reflect the latest generated numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 9958
$ws1.Range("F16").Value = 2031
$ws1.Range("F20").Value = 1577
$ws1.Range("F22").Value = 44
$ws1.Range("F23").Value = 217
$ws1.Range("F28").Value = 349
$ws1.Range("F33").Value = 282
$ws1.Range("F36").Value = 395
$ws1.Range("F38").Value = 424

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 33

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9958
$ws4.Range("F22").Value = 1577
$ws4.Range("F24").Value = 44
$ws4.Range("F25").Value = 217
$ws4.Range("F30").Value = 349
$ws4.Range("F32").Value = 33
$ws4.Range("F39").Value = 282
$ws4.Range("F43").Value = 395
$ws4.Range("F45").Value = 424
